# "traigo los precios de todos los bonos" -- the sheet used to hold a single
# bond's fechaHora/ultimoPrecio pair (merge output from one pandas merge).
# Now it holds the merged prices for every bond: fechaHora plus three
# ultimoPrecio_x/ultimoPrecio_y pairs (one per successive pandas merge) and
# a final plain ultimoPrecio column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---------------------------------------------------
# B1 "fechaHora" and C1 "ultimoPrecio" already exist. Add the new bond price
# columns D1:I1, then turn C1 into the first "_x" column and finish with a
# plain "ultimoPrecio" in I1.
$ws.Range("D1").Value = "ultimoPrecio_y"
$ws.Range("E1").Value = "ultimoPrecio_x"
$ws.Range("F1").Value = "ultimoPrecio_y"
$ws.Range("G1").Value = "ultimoPrecio_x"
$ws.Range("H1").Value = "ultimoPrecio_y"
$ws.Range("I1").Value = "ultimoPrecio"

# Give the new header cells the same look (bold, centered, bordered) as the
# existing header cells by copying C1's format onto them.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1:I1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# C1 itself becomes "ultimoPrecio_x" (was "ultimoPrecio").
$ws.Range("C1").Value = "ultimoPrecio_x"

# --- Data row (row 2) ------------------------------------------------------
$ws.Range("B2").Value = "20/08/2021"
$ws.Range("C2").Value = 36.21
$ws.Range("D2").Value = 39
$ws.Range("E2").Value = 35.27
$ws.Range("F2").Value = 37.99
$ws.Range("G2").Value = 32.7
$ws.Range("H2").Value = 37.2
$ws.Range("I2").Value = 36.6
